$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Build the two new border styles once (top+bottom, and top+right+bottom)
# on sheet1's C1/D1, then fan them out by copying the format so the style
# table only ends up with the two extra shared style entries that the
# target workbook expects (instead of rebuilding borders edge-by-edge on
# every cell, which would leave unused intermediate style entries behind).

$c1s1 = $ws1.Range("C1")
$c1s1.ClearFormats()
$c1s1.Borders.Item(8).LineStyle = 1   # top
$c1s1.Borders.Item(9).LineStyle = 1   # bottom

$d1s1 = $ws1.Range("D1")
$d1s1.ClearFormats()
$d1s1.Borders.Item(8).LineStyle = 1   # top
$d1s1.Borders.Item(10).LineStyle = 1  # right
$d1s1.Borders.Item(9).LineStyle = 1   # bottom

$c1s1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)

$d1s1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Text updates: anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty inline-string cell
$ws2.Range("G5").ClearContents()
